# Applies the "add stat desc, plots and readme" commit:
#  - removes the (empty) "c_netneh" sheet
#  - updates the data-type labels on the "Feuil1" codebook sheet
#    (int -> int32/int16, category -> categorie)
#  - adds a new "Modifications" column (F) documenting a data fix
#  - tweaks a few column widths / selections to match the saved workbook

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Drop the empty "c_netneh" worksheet entirely.
# ---------------------------------------------------------------------------
$wsNetneh = $wb.Worksheets.Item("c_netneh")
$wsNetneh.Delete()

# ---------------------------------------------------------------------------
# 2. "Feuil1" codebook sheet: refresh dtype labels + add Modifications column
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Feuil1")

# New header cell for the extra column.
$ws1.Range("F1").Value = "Modifications"

# Refresh the "type de donnée" column with the more precise numpy dtypes.
$ws1.Range("B2").Value = "int32"   # ident
$ws1.Range("B3").Value = "int32"   # annee
$ws1.Range("B4").Value = "int16"   # trimestre
$ws1.Range("B5").Value = "categorie"  # statut
$ws1.Range("B6").Value = "categorie"  # qualite
$ws1.Range("B7").Value = "categorie"  # etat
$ws1.Range("B11").Value = "int32"  # ib_

# Document the NaN -> -1 fix applied to ib_ (Indice brut).
$ws1.Range("F11").Value = "remplace Nan par -1"

# Column widths / dimension follow from the edits above; just tune E & F.
$ws1.Columns.Item(5).ColumnWidth = 46.5
$ws1.Columns.Item(6).ColumnWidth = 21.666666666666668

# Selection ends up on the new column F (whole-column selection).
$ws1.Range("F1:F1048576").Select()

# ---------------------------------------------------------------------------
# 3. "etat" sheet: just a saved-selection change.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("etat")
$ws2.Range("B14").Select()

# ---------------------------------------------------------------------------
# 4. "qualite" sheet: selection + slightly wider column B.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("qualite")
$ws3.Columns.Item(2).ColumnWidth = 15.666666666666666
$ws3.Range("A2:A11").Select()

# ---------------------------------------------------------------------------
# Leave the workbook focused back on "Feuil1", matching its original state.
# ---------------------------------------------------------------------------
$ws1.Activate()
